# Add a couple of new publications to the single-cell RNA-seq comparison workbook.

$wb = $excel.ActiveWorkbook

$wsHuman   = $wb.Worksheets.Item(1)   # "Human"
$wsMouse   = $wb.Worksheets.Item(2)   # "Mouse"
$wsOrganoid = $wb.Worksheets.Item(3)  # "Human organoid"

# ---------------------------------------------------------------------------
# 1) Human sheet: append a new row (row 10) for the Onorati et al. publication
# ---------------------------------------------------------------------------
$wsHuman.Range("A10").Value = '<a href="https://www.nature.com/articles/s41422-018-0053-3"  target="_blank">Onorati</a>'
$wsHuman.Range("B10").Value = "C1"
$wsHuman.Range("C10").Value = "Full-length"
$wsHuman.Range("D10").Value = "5-20pcw"

# Widen column A on the Human sheet to fit the new, longer reference text.
$wsHuman.Columns.Item(1).ColumnWidth = 86

# ---------------------------------------------------------------------------
# 2) Human organoid sheet: append a new row (row 5) that repeats the
#    reference (Madhavan et al.) already used on row 4.
# ---------------------------------------------------------------------------
$wsOrganoid.Range("A5").Value = '<a href="https://www.nature.com/articles/s41592-018-0081-4" target=”_blank”>Madhavan</a>'

# ---------------------------------------------------------------------------
# 3) Update the view state: selections on each sheet, and make the
#    "Human organoid" sheet the active/selected tab.
# ---------------------------------------------------------------------------
$wsHuman.Activate()
$wsHuman.Range("A10").Select()

$wsMouse.Activate()
$wsMouse.Range("A11").Select()

$wsOrganoid.Activate()
$wsOrganoid.Range("A5").Select()
